$d = $word.ActiveDocument

# The document starts with:
#   1: "Orphan Report" (Title)
#   2: "" (empty paragraph)
#   3..17: alternating PUMP:RISK:* tag / requirement text / "is an orphan tag" lines
#
# Target: drop the empty paragraph and collapse every PUMP:RISK:* related
# paragraph into a single bold summary paragraph right after the title.

# 1) Remove the blank paragraph that follows the title.
$d.Paragraphs.Item(2).Range.Delete()

# 2) Turn the (now second) paragraph into the new summary text.
$newText = "These are the orphan tags that were found in the documents: "
$summaryPara = $d.Paragraphs.Item(2)
$summaryPara.Range.Text = $newText

# 3) Delete every paragraph that follows the summary paragraph - these held
#    the old PUMP:RISK tag / requirement / orphan-tag-notice text.
$paraCount = $d.Paragraphs.Count
if ($paraCount -gt 2) {
    $deleteStart = $d.Paragraphs.Item(3).Range.Start
    $deleteEnd = $d.Paragraphs.Item($paraCount).Range.End
    $d.Range($deleteStart, $deleteEnd).Delete()
}

# 4) Bold just the new text (not the paragraph mark) so only a run-level
#    <w:rPr><w:b/></w:rPr> is produced, matching the target markup.
$textRange = $d.Range($summaryPara.Range.Start, $summaryPara.Range.Start + $newText.Length)
$textRange.Bold = 1
